$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cyclic rotation of the date/volume/price data across rows 2,6,4,11,9,7,10,8
# (rows 3 and 5 are unchanged). Values below are the new targets per row,
# matching the OOXML diff.

$ws.Range("D2").Value = 44761
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("S2").Value = 1025

$ws.Range("D4").Value = 44533
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("S4").Value = 825

$ws.Range("D6").Value = 44798
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 1075

$ws.Range("D7").Value = 44890
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22250
$ws.Range("S7").Value = 1112

$ws.Range("D8").Value = 44357
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("S8").Value = 725

$ws.Range("D9").Value = 45092
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 24000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 24333
$ws.Range("S9").Value = 1217

$ws.Range("D10").Value = 44708
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 21000
$ws.Range("P10").Value = 20500
$ws.Range("S10").Value = 1025

$ws.Range("D11").Value = 44320
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 16500
$ws.Range("S11").Value = 825
